$wb = $excel.ActiveWorkbook

# "revision" sheet holds the revision history table in column A (REV N -> V.N)
$revSheet = $wb.Worksheets.Item("revision")

$revSheet.Range("A1").Value = "V.0"
$revSheet.Range("A2").Value = "V.1"
$revSheet.Range("A3").Value = "V.2"
$revSheet.Range("A4").Value = "V.3"
$revSheet.Range("A5").Value = "V.4"
$revSheet.Range("A6").Value = "V.5"
$revSheet.Range("A7").Value = "V.6"
$revSheet.Range("A8").Value = "V.7"
$revSheet.Range("A9").Value = "V.8"
$revSheet.Range("A10").Value = "V.9"
$revSheet.Range("A11").Value = "V.x"

# Update timestamp on the last revision row, and narrow column A
$revSheet.Range("B11").Value = 43167.4860474883
# The ColumnWidth property adds a fixed ~0.8333 char padding when it is
# serialized back to the OOXML <col width> attribute, so subtract that
# offset here in order to land on an exact width of 6 in the saved file.
$revSheet.Columns.Item(1).ColumnWidth = 5.166666666666667

# "Parts - Consoles" sheet holds a mirrored timestamp in BA2
$dataSheet = $wb.Worksheets.Item("Parts - Consoles")
$dataSheet.Range("BA2").Value = 43167.48612357685
